$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("挑战组")
$ws2 = $wb.Worksheets.Item("中坚组")
$ws3 = $wb.Worksheets.Item("传奇组")

$newLink = "https://www.bilibili.com/video/BV1Pg411r7V5"

# Fill the new video link text down column H (sheet1, rows 2-17)
for ($r = 2; $r -le 17; $r++) {
    $ws1.Cells.Item($r, 8).Value = $newLink
}

# Fill the new video link text down column F (sheet2, rows 2-17)
for ($r = 2; $r -le 17; $r++) {
    $ws2.Cells.Item($r, 6).Value = $newLink
}

# Fill the new video link text down column F (sheet3, rows 2-16)
for ($r = 2; $r -le 16; $r++) {
    $ws3.Cells.Item($r, 6).Value = $newLink
}

# sheet3's F17 already carried the (old) link as a real hyperlink; remove the
# hyperlink object and its special formatting, replacing it with the new
# plain-text link value so it matches the rest of the filled-down column.
$ws3.Hyperlinks.Delete()
$ws3.Cells.Item(17, 6).Value = $newLink
$ws3.Cells.Item(17, 6).Style = "常规"
$wb.Styles.Item("超链接").Delete()

# Restore print setup on the legend-group sheet
$ws3.PageSetup.PaperSize = 9
$ws3.PageSetup.Orientation = 1

# Reproduce the selections/active sheet left behind by the edit
$ws2.Range("F2:F17").Select()
$ws3.Range("B41").Select()
$ws1.Range("H2:H17").Select()
